$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header cell: copy style from E1 (bold, border, centered) then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# F2:F63 data cells: plain "time_taken" timestamps (no special style)
$ws.Range("F2").Value = "2021-10-05 13:41:19.874684"
$ws.Range("F3").Value = "2021-10-05 13:41:19.874695"
$ws.Range("F4").Value = "2021-10-05 13:41:19.874699"
$ws.Range("F5").Value = "2021-10-05 13:41:19.874701"
$ws.Range("F6").Value = "2021-10-05 13:41:19.874704"
$ws.Range("F7").Value = "2021-10-05 13:41:19.874707"
$ws.Range("F8").Value = "2021-10-05 13:41:19.874709"
$ws.Range("F9").Value = "2021-10-05 13:41:19.874712"
$ws.Range("F10").Value = "2021-10-05 13:41:19.874714"
$ws.Range("F11").Value = "2021-10-05 13:41:19.874717"
$ws.Range("F12").Value = "2021-10-05 13:41:19.874719"
$ws.Range("F13").Value = "2021-10-05 13:41:19.874722"
$ws.Range("F14").Value = "2021-10-05 13:41:19.874724"
$ws.Range("F15").Value = "2021-10-05 13:41:19.874727"
$ws.Range("F16").Value = "2021-10-05 13:41:19.874729"
$ws.Range("F17").Value = "2021-10-05 13:41:19.874732"
$ws.Range("F18").Value = "2021-10-05 13:41:19.874734"
$ws.Range("F19").Value = "2021-10-05 13:41:19.874737"
$ws.Range("F20").Value = "2021-10-05 13:41:19.874739"
$ws.Range("F21").Value = "2021-10-05 13:41:19.874742"
$ws.Range("F22").Value = "2021-10-05 13:41:19.874744"
$ws.Range("F23").Value = "2021-10-05 13:41:19.874747"
$ws.Range("F24").Value = "2021-10-05 13:41:19.874749"
$ws.Range("F25").Value = "2021-10-05 13:41:19.874752"
$ws.Range("F26").Value = "2021-10-05 13:41:19.874755"
$ws.Range("F27").Value = "2021-10-05 13:41:19.874757"
$ws.Range("F28").Value = "2021-10-05 13:41:19.874760"
$ws.Range("F29").Value = "2021-10-05 13:41:19.874762"
$ws.Range("F30").Value = "2021-10-05 13:41:19.874764"
$ws.Range("F31").Value = "2021-10-05 13:41:19.874767"
$ws.Range("F32").Value = "2021-10-05 13:41:19.874769"
$ws.Range("F33").Value = "2021-10-05 13:41:19.874772"
$ws.Range("F34").Value = "2021-10-05 13:41:19.874775"
$ws.Range("F35").Value = "2021-10-05 13:41:19.874777"
$ws.Range("F36").Value = "2021-10-05 13:41:19.874780"
$ws.Range("F37").Value = "2021-10-05 13:41:19.874782"
$ws.Range("F38").Value = "2021-10-05 13:41:19.874785"
$ws.Range("F39").Value = "2021-10-05 13:41:19.874787"
$ws.Range("F40").Value = "2021-10-05 13:41:19.874790"
$ws.Range("F41").Value = "2021-10-05 13:41:19.874792"
$ws.Range("F42").Value = "2021-10-05 13:41:19.874795"
$ws.Range("F43").Value = "2021-10-05 13:41:19.874798"
$ws.Range("F44").Value = "2021-10-05 13:41:19.874800"
$ws.Range("F45").Value = "2021-10-05 13:41:19.874803"
$ws.Range("F46").Value = "2021-10-05 13:41:19.874805"
$ws.Range("F47").Value = "2021-10-05 13:41:19.874807"
$ws.Range("F48").Value = "2021-10-05 13:41:19.874810"
$ws.Range("F49").Value = "2021-10-05 13:41:19.874812"
$ws.Range("F50").Value = "2021-10-05 13:41:19.874815"
$ws.Range("F51").Value = "2021-10-05 13:41:19.874817"
$ws.Range("F52").Value = "2021-10-05 13:41:19.874820"
$ws.Range("F53").Value = "2021-10-05 13:41:19.874822"
$ws.Range("F54").Value = "2021-10-05 13:41:19.874825"
$ws.Range("F55").Value = "2021-10-05 13:41:19.874827"
$ws.Range("F56").Value = "2021-10-05 13:41:19.874830"
$ws.Range("F57").Value = "2021-10-05 13:41:19.874832"
$ws.Range("F58").Value = "2021-10-05 13:41:19.874835"
$ws.Range("F59").Value = "2021-10-05 13:41:19.874837"
$ws.Range("F60").Value = "2021-10-05 13:41:19.874840"
$ws.Range("F61").Value = "2021-10-05 13:41:19.874842"
$ws.Range("F62").Value = "2021-10-05 13:41:19.874845"
$ws.Range("F63").Value = "2021-10-05 13:41:19.874847"
